$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> changes (date text, and updated D/E/G/H numeric values where applicable)
$rows = @(
    @{ Row = 3;  Date = "28-07-2022"; D = 1; G = 1 },
    @{ Row = 4;  Date = "01-08-2022"; D = 1; E = 1; H = 0 },
    @{ Row = 5;  Date = "04-08-2022"; D = 1; E = 1; H = 0 },
    @{ Row = 6;  Date = "08-08-2022" },
    @{ Row = 7;  Date = "11-08-2022"; D = 1; E = 1; H = 0 },
    @{ Row = 8;  Date = "15-08-2022" },
    @{ Row = 9;  Date = "18-08-2022" },
    @{ Row = 10; Date = "22-08-2022"; D = 1; E = 1; H = 0 },
    @{ Row = 11; Date = "25-08-2022" },
    @{ Row = 12; Date = "29-08-2022"; D = 1; E = 1; H = 0 },
    @{ Row = 13; Date = "01-09-2022" },
    @{ Row = 14; Date = "05-09-2022" },
    @{ Row = 15; Date = "08-09-2022" },
    @{ Row = 16; Date = "12-09-2022" },
    @{ Row = 17; Date = "15-09-2022" },
    @{ Row = 18; Date = "19-09-2022" },
    @{ Row = 19; Date = "22-09-2022" },
    @{ Row = 20; Date = "26-09-2022" },
    @{ Row = 21; Date = "29-09-2022" }
)

foreach ($item in $rows) {
    $r = $item.Row
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $item.Date
    if ($item.ContainsKey("D")) { $ws.Cells.Item($r, 4).Value = $item.D }
    if ($item.ContainsKey("E")) { $ws.Cells.Item($r, 5).Value = $item.E }
    if ($item.ContainsKey("G")) { $ws.Cells.Item($r, 7).Value = $item.G }
    if ($item.ContainsKey("H")) { $ws.Cells.Item($r, 8).Value = $item.H }
}
